$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: "Year" -> "Dates"
$ws.Range("A1").Value = "Dates"

# Dates (column A, as plain text) and values (column B) for rows 2-11
$dates = @(
    "2023-01-01",
    "2024-01-01",
    "2025-01-01",
    "2026-01-01",
    "2027-01-01",
    "2028-01-01",
    "2029-01-01",
    "2030-01-01",
    "2031-01-01",
    "2032-01-01"
)

$values = @(-0.22, -1.88, -1.65, -1.72, -3.37, -3.77, -0.01, -4.07, -1.4, -2.99)

$lastRow = 1 + $dates.Length

# Clear the old date-serial number formatting up front so the new text
# values don't inherit/allocate per-cell number-format styles.
$ws.Range("A2:A$lastRow").ClearFormats()

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    # Leading apostrophe forces the value to be stored as literal text
    # instead of being auto-parsed back into a date serial number.
    $ws.Range("A$row").Formula = "'" + $dates[$i]
    $ws.Range("B$row").Value = $values[$i]
}

# Re-clear formatting so column A ends up with the default (unstyled) cell
# style rather than the transient quote-prefixed style used while typing.
$ws.Range("A2:A$lastRow").ClearFormats()
